$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.458.46'
$ws.Range("E2").Value = '  +1.75%  '

$ws.Range("D3").Value = '2.924.06'
$ws.Range("E3").Value = '  +4.90%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '352.84'
$ws.Range("E5").Value = '  -0.42%  '

$ws.Range("D6").Value = '112.95'
$ws.Range("E6").Value = '  +3.86%  '

$ws.Range("D7").Value = '0.559'
$ws.Range("E7").Value = '  +0.68%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").Value = '0.627'
$ws.Range("E9").Value = '  +0.44%  '

$ws.Range("D10").Value = '40.32'
$ws.Range("E10").Value = '  +1.12%  '

$ws.Range("D11").Value = '0.0862'
$ws.Range("E11").Value = '  +2.99%  '

$ws.Range("E12").Value = '  +0.59%  '

$ws.Range("D13").Value = '20.17'
$ws.Range("E13").Value = '  +1.06%  '

$ws.Range("D14").Value = '7.85'
$ws.Range("E14").Value = '  +1.42%  '

$ws.Range("D15").Value = '3.384.05'
$ws.Range("E15").Value = '  +4.73%  '

$ws.Range("D16").Value = '2.927.72'
$ws.Range("E16").Value = '  +4.25%  '

$ws.Range("D17").Value = '0.998'
$ws.Range("E17").Value = '  +6.47%  '

$ws.Range("D18").Value = '52.457.04'
$ws.Range("E18").Value = '  +1.84%  '

$ws.Range("D19").Value = '7.75'
$ws.Range("E19").Value = '  +0.00%  '

$ws.Range("D20").Value = '3.36'
$ws.Range("E20").Value = '  +6.23%  '

$ws.Range("E21").Value = '  +7.78%  '

$ws.Range("E22").Value = '  +1.43%  '

$ws.Range("D23").Value = '71.25'
$ws.Range("E23").Value = '  +1.46%  '

$ws.Range("D24").Value = '271.79'

$ws.Range("E25").Value = '  +1.61%  '

$ws.Range("E26").Value = '  +3.70%  '

$ws.Range("E27").Value = '  +0.15%  '

$ws.Range("E28").Value = '  -0.58%  '

$ws.Range("E29").Value = '  +2.67%  '

$ws.Range("D30").Value = '38.42'
$ws.Range("E30").Value = '  +3.62%  '

$ws.Range("E31").Value = '  +1.13%  '

$ws.Range("B32").Value = 'RenderToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D32").Value = '6.19'
$ws.Range("E32").Value = '  +9.06%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '6.45'
$ws.Range("E33").Value = '  +3.58%  '

$ws.Range("D34").Value = '53.07'
$ws.Range("E34").Value = '  +2.29%  '

$ws.Range("D35").Value = '0.0938'
$ws.Range("E35").Value = '  +10.09%  '

$ws.Range("E36").Value = '  +3.37%  '

$ws.Range("E37").Value = '  -0.25%  '

$ws.Range("D38").Value = '3.35'
$ws.Range("E38").Value = '  +6.94%  '

$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").Value = '18.90'
$ws.Range("E39").Value = '  +0.68%  '

$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '2.09'
$ws.Range("E40").Value = '  +5.68%  '

$ws.Range("D41").Value = '2.73'
$ws.Range("E41").Value = '  +9.70%  '

$ws.Range("D42").Value = '24.67'
$ws.Range("E42").Value = '  +12.94%  '

$ws.Range("E43").Value = '  +2.19%  '

$ws.Range("D44").Value = '122.23'
$ws.Range("E44").Value = '  +2.62%  '

$ws.Range("E45").Value = '  +0.80%  '

$ws.Range("D46").Value = '3.57'
$ws.Range("E46").Value = '  +5.34%  '

$ws.Range("D47").Value = '2.216.84'
$ws.Range("E47").Value = '  +4.62%  '

$ws.Range("E48").Value = '  +6.56%  '

$ws.Range("E49").Value = '  +23.87%  '

$ws.Range("D50").Value = '0.954'
$ws.Range("E50").Value = '  +4.70%  '

$ws.Range("D51").Value = '0.0328'
$ws.Range("E51").Value = '  +13.62%  '
